# Update student 48966's grade row (row 9) on the active worksheet:
# B9: 76 -> 78 (Calculated Grade)
# C9: 58.66 -> 62.82853277352308 (Weighted Grade)
# D9: -17.34 -> -15.17146722647692 (Difference)
# F9: 17.34 -> 15.17146722647692 (Penalty)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 78
$ws.Range("C9").Value = 62.82853277352308
$ws.Range("D9").Value = -15.17146722647692
$ws.Range("F9").Value = 15.17146722647692
